$wb = $excel.ActiveWorkbook

$wsTrans = $wb.Worksheets.Item("Transformations")
$wsViews = $wb.Worksheets.Item("Views")

# Populate new shared strings in the same order they first appear
# in the saved workbook (test_doub, test2, doubling_time(...)):
# --- Transformations sheet: add new row 18 (A18) ---
$wsTrans.Range("A18").Value = "test_doub"

# --- Views sheet: add new row 5 (A5, B5) ---
$wsViews.Range("A5").Value = "test2"
$wsViews.Range("B5").Value = "test_doub"

# --- Transformations sheet: add new row 18 (B18) ---
$wsTrans.Range("B18").Value = "doubling_time(combine_samples_od, plate_01_time.OD;max_od=0.2)"

# Update selection on Views sheet to B6 (matches the post-edit diff)
$wsViews.Range("B6").Select()

# Make sure Transformations tab stays the active/selected sheet
$wsTrans.Select()
$wsTrans.Range("B18").Select()
